$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2 = @(0.3146293670956112, 0.003648818996706889)
    3 = @(0.3067188829323158, 0.006251139570123933)
    4 = @(0.3566314791467504, 0.007518210472485615)
    5 = @(0.3302645193144947, 0.005578099084910562)
    6 = @(0.305696041515249, 0.005087245722568943)
    7 = @(0.2614210394248254, 0.003239970372751299)
    8 = @(0.3025824831866997, 0.004538753918966533)
    9 = @(0.3021027523813148, 0.005217041126488467)
    10 = @(0.3386647529303305, 0.005615282256397132)
    11 = @(0.2587795962665891, 0.003925219762725651)
    12 = @(0.3240203674843917, 0.007021938290808478)
    13 = @(0.2914894387684464, 0.00631817727714739)
    14 = @(0.2067006618075515, 0.002383590027845342)
    15 = @(0.1264794623308638, 0.005508080263651596)
    16 = @(0.09151841352835277, 0.005697387575726859)
    17 = @(0.09061337401526949, 0.005532012816701792)
    18 = @(0.04156291724602371, 0.001320968378429379)
}

foreach ($row in $values.Keys) {
    $bc = $values[$row]
    $ws.Cells.Item($row, 2).Value = $bc[0]
    $ws.Cells.Item($row, 3).Value = $bc[1]
}